# Adds the 28.9.2018 work-log entry (row 13) to the "Eetu Pihamäki" sheet.
# The hour totals on both sheets are driven by formulas, so they recalculate
# automatically once the new row's time values are entered.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Eetu Pihamäki")

# Pvm (date), Aloitusklo (start time), Lopetusklo (end time), Sprint, Tehtävä
$ws.Range("A13").Value = 43371
$ws.Range("B13").Value = 0.375
$ws.Range("C13").Value = 0.63194444444444442
$ws.Range("E13").Value = 1
$ws.Range("F13").Value = "6h etsin 9 vertailuvaatimukseen tietoa IdM-järjestelmistä. https://github.com/Eetu95/Open-source-IdM-solution/blob/master/Eetun%20muistiinpanoja/Ty%C3%B6t%20-%2028.9.2018.txt"

# The long "Tehtävä" text wraps onto multiple lines, matching the other
# multi-line rows in this table (e.g. row 9 / row 12 are also 60pt tall).
$ws.Rows(13).RowHeight = 60

$excel.Calculate()

Write-Host "Added row 13; Eetu Pihamäki total (C5):" $ws.Range("C5").Value()
